# Update the "Fruta, Femacal de La Calera - Coco" weekly price sheet.
# The Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns are updated
# for rows 2-20 to reflect the new weekly data (rows effectively got
# their date/price/volume data reshuffled to new dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedioPonderado, PrecioPorKg)
$rows = @{
    2  = @(44413, 45, 20000, 20000, 20000, 1000)
    3  = @(44377, 25, 20000, 20000, 20000, 1000)
    5  = @(44300, 45, 22000, 22000, 22000, 1100)
    6  = @(44389, 20, 20000, 20000, 20000, 1000)
    7  = @(44298, 65, 22000, 22000, 22000, 1100)
    8  = @(44385, 36, 20000, 20000, 20000, 1000)
    9  = @(44403, 50, 20000, 20000, 20000, 1000)
    11 = @(44301, 38, 22000, 22000, 22000, 1100)
    12 = @(44448, 30, 22000, 22000, 22000, 1100)
    14 = @(44292, 30, 25000, 25000, 25000, 1250)
    15 = @(44400, 45, 20000, 20000, 20000, 1000)
    16 = @(44376, 38, 20000, 20000, 20000, 1000)
    17 = @(44305, 20, 22000, 22000, 22000, 1100)
    18 = @(44382, 24, 20000, 20000, 20000, 1000)
    19 = @(44307, 30, 22000, 22000, 22000, 1100)
    20 = @(44291, 70, 25000, 25000, 25000, 1250)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]  # D: Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]  # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]  # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]  # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]  # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]  # S: Precio $/Kg
}
